$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (prices/volume %s updated by the scrape run).
$updates = @{
    "D2" = "64.343.44"
    "E2" = "  +0.08%  "
    "D3" = "3.496.18"
    "E3" = "  +0.12%  "
    "E4" = "  +0.04%  "
    "D5" = "589.06"
    "E5" = "  +0.42%  "
    "D6" = "134.20"
    "E8" = "  +0.55%  "
    "D9" = "7.61"
    "E9" = "  +5.95%  "
    "E10" = "  +0.43%  "
    "E11" = "  +3.63%  "
    "D12" = "4.093.00"
    "E12" = "  +0.21%  "
    "E13" = "  +0.67%  "
    "E14" = "  +0.26%  "
    "D15" = "3.496.93"
    "E15" = "  +0.22%  "
    "D16" = "64.343.54"
    "D17" = "25.44"
    "E17" = "  +1.24%  "
    "D18" = "10.01"
    "E19" = "  +0.85%  "
    "D20" = "13.53"
    "E20" = "  -0.56%  "
    "D21" = "388.35"
    "E21" = "  +0.21%  "
    "D22" = "0.581"
    "E22" = "  +2.83%  "
    "D23" = "3.635.98"
    "E23" = "  +0.13%  "
    "D24" = "74.27"
    "E24" = "  -0.44%  "
    "E25" = "  +0.05%  "
    "E26" = "  -1.09%  "
    "D27" = "0.0000115"
    "E27" = "  +2.41%  "
    "E28" = "  -0.10%  "
    "D29" = "7.38"
    "E29" = "  -0.03%  "
    "E30" = "  +1.37%  "
    "D31" = "1.50"
    "E31" = "  -3.30%  "
    "E32" = "  -1.16%  "
    "E33" = "  +5.39%  "
    "D34" = "3.522.60"
    "E34" = "  +0.29%  "
    "D36" = "23.36"
    "E36" = "  -0.45%  "
    "D37" = "5.33"
    "E37" = "  +1.53%  "
    "D38" = "6.93"
    "E38" = "  +1.42%  "
    "D39" = "1.54"
    "E39" = "  +0.78%  "
    "D40" = "165.55"
    "E40" = "  +2.48%  "
    "D41" = "0.0787"
    "E41" = "  +0.97%  "
    "D42" = "0.808"
    "E42" = "  +0.50%  "
    "E43" = "  +0.03%  "
    "E44" = "  +0.69%  "
    "D45" = "24.86"
    "E45" = "  -1.60%  "
    "D46" = "1.18"
    "E46" = "  +0.21%  "
    "E47" = "  -0.16%  "
    "D48" = "6.82"
    "E48" = "  +1.27%  "
    "E49" = "  +2.78%  "
    "D50" = "2.402.10"
    "E50" = "  -2.85%  "
    "E51" = "  -0.02%  "
}

foreach ($addr in $updates.Keys) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    # Force text so numeric-looking strings (e.g. "134.20", "7.61") are not
    # coerced into numbers (which would drop trailing zeros / change type),
    # matching the source data which stores these as inline strings.
    $r.NumberFormat = "@"
    $r.Value = $updates[$addr]
    $r.Style = $origStyle
}
